# Swap the data of row 11 and row 12 (columns B..AC), keeping row labels
# in column A untouched. Columns A, C, D, E, J happen to hold identical
# values in both rows, so swapping the full B:AC range reproduces the
# diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 11
$row2 = 12

$rng1 = $ws.Range("B$row1`:AC$row1")
$rng2 = $ws.Range("B$row2`:AC$row2")

$vals1 = $rng1.Value2
$vals2 = $rng2.Value2

$rng1.Value2 = $vals2
$rng2.Value2 = $vals1
